# Insert a new row at position 14 (pushes existing rows 14..116 down to 15..117)
# and populate it with a new weekly price record for Berenjena (matches the
# commit "Fruta / hortaliza, semanal" weekly data-refresh pattern).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = "2023-08-31"
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = 100112001
$ws.Range("G14").Value = "Berenjena"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("N14").Value = "$/caja 60 unidades"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 167
$ws.Range("Q14").Value = 60
$ws.Range("R14").Value = "Hortaliza"
